# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" sheet (fund holdings detail) right after "2021-Q4"
#   and before "总计" (matching the per-quarter sheet layout used by
#   2020-Q4 .. 2021-Q4).
# - Rebuilds the "总计" (totals) sheet with a new first data row for
#   2022-Q1 (5 funds, 1.9 亿元), pushing the previously existing rows down
#   by one.

$wb = $excel.ActiveWorkbook

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 0. Grab a style "donor" cell (s=2: bold, thin border, center/top align)
#    from the existing "2021-Q4" sheet (same look the header + index
#    cells use everywhere in this workbook). We read it from a sheet we
#    never delete, since deleting a sheet invalidates a pending
#    Copy()/clipboard reference to it.
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("B1").Copy()

$totalSheetOld = $wb.Worksheets.Item("总计")

# Remove the old "总计" sheet - it will be rebuilt from scratch after the
# new "2022-Q1" sheet so both end up with the right sheetId / position.
$totalSheetOld.Delete()

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet (fund holding detail), placed right after
#    "2021-Q4".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4Sheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"
$q1.Range("B1:H1").PasteSpecial($fmt)

$q1Data = @(
    @(0, "009993", "嘉实前沿创新混合",             "16.96", "88.41", "4.97", "0.8429", 6),
    @(1, "004450", "嘉实前沿科技沪港深股票",       "21.72", "88.32", "3.75", "0.8145", 10),
    @(2, "011851", "天弘先进制造混合型证券投资基金A", "2.72", "91.41", "5.31", "0.1444", 6),
    @(3, "003292", "嘉实优势成长灵活配置混合",     "1.04",  "92.17", "5.44", "0.0566", 8),
    @(4, "011852", "天弘先进制造混合型证券投资基金C", "0.70", "91.41", "5.31", "0.0372", 6)
)

foreach ($row in $q1Data) {
    $r = [int]$row[0] + 2

    $q1.Cells.Item($r, 1).Value = [int]$row[0]
    $q1.Cells.Item($r, 1).PasteSpecial($fmt)

    $q1.Range($q1.Cells.Item($r, 2), $q1.Cells.Item($r, 7)).NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Rebuild "总计" sheet right after "2022-Q1", with the new 2022-Q1
#    row inserted at the top of the data and everything else shifted
#    down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$total.Range("B1:D1").PasteSpecial($fmt)

$totalData = @(
    @(0, "2022-Q1", 5, 1.9),
    @(1, "2021-Q4", 9, 3.74),
    @(2, "2021-Q3", 6, 1.03),
    @(3, "2021-Q2", 1, 0.04),
    @(4, "2021-Q1", 2, 0.06),
    @(5, "2020-Q4", 2, 0.09)
)

foreach ($row in $totalData) {
    $r = [int]$row[0] + 2

    $total.Cells.Item($r, 1).Value = [int]$row[0]
    $total.Cells.Item($r, 1).PasteSpecial($fmt)

    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

$q4Sheet.Activate()
